# Insert a new record row before row 597 on the active sheet.
# This pushes the existing rows 597-641 down to 598-642 (Excel preserves
# their data/formatting automatically), and we then fill in the new
# row 597 with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 597 (shifts 597:641 -> 598:642)
$ws.Rows.Item(597).Insert()

# Populate the new row 597 with the new data record
$ws.Range("A597").Value = 8
$ws.Range("B597").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C597").Value = 'Coquimbo'
$ws.Range("D597").Value = 44769
$ws.Range("E597").Value = 4
$ws.Range("F597").Value = 100112024
$ws.Range("G597").Value = 'Choclo'
$ws.Range("H597").Value = 'Dulce o Americano'
$ws.Range("I597").Value = 'Primera'
$ws.Range("J597").Value = 500
$ws.Range("K597").Value = 40000
$ws.Range("L597").Value = 41000
$ws.Range("M597").Value = 40500
$ws.Range("N597").Value = '$/malla 70 unidades'
$ws.Range("O597").Value = 'Región de Arica y Parinacota'
$ws.Range("P597").Value = 579
$ws.Range("Q597").Value = 70
$ws.Range("R597").Value = 'Hortaliza'

# Keep the date column's number format consistent with the rest of column D
$ws.Range("D597").NumberFormat = $ws.Range("D598").NumberFormat
